$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 22:21"

# --- Swap country-name labels that changed order in the source list ---
# (the underlying per-row statistics below already correspond to the
#  country that now occupies that row after the swap)
$tmpA = $ws.Range("A131").Value2
$tmpB = $ws.Range("A132").Value2
$ws.Range("A131").Value = $tmpB
$ws.Range("A132").Value = $tmpA
$tmpA = $ws.Range("A148").Value2
$tmpB = $ws.Range("A149").Value2
$ws.Range("A148").Value = $tmpB
$ws.Range("A149").Value = $tmpA
$tmpA = $ws.Range("A202").Value2
$tmpB = $ws.Range("A203").Value2
$ws.Range("A202").Value = $tmpB
$ws.Range("A203").Value = $tmpA
$tmpA = $ws.Range("A210").Value2
$tmpB = $ws.Range("A211").Value2
$ws.Range("A210").Value = $tmpB
$ws.Range("A211").Value = $tmpA

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 4150887
$ws.Range("C4").Value = 50012
$ws.Range("D4").Value = 1961011
$ws.Range("E4").Value = 2042901
$ws.Range("G4").Value = 792
$ws.Range("H4").Value = 146975
$ws.Range("B6").Value = 1288130
$ws.Range("C6").Value = 48446
$ws.Range("D6").Value = 817593
$ws.Range("E6").Value = 439892
$ws.Range("B8").Value = 408052
$ws.Range("C8").Value = 13104
$ws.Range("D8").Value = 236260
$ws.Range("E8").Value = 165699
$ws.Range("G8").Value = 153
$ws.Range("H8").Value = 6093
$ws.Range("B21").Value = 205142
$ws.Range("C21").Value = 672
$ws.Range("E21").Value = 6955
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 9187
$ws.Range("B30").Value = 78763
$ws.Range("C30").Value = 53
$ws.Range("G30").Value = 10
$ws.Range("H30").Value = 5676
$ws.Range("B31").Value = 78148
$ws.Range("C31").Value = 891
$ws.Range("D31").Value = 33455
$ws.Range("E31").Value = 39254
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 5439
$ws.Range("E51").Value = 3610
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 133
$ws.Range("D69").Value = 7135
$ws.Range("E69").Value = 8203
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 263
$ws.Range("B77").Value = 13129
$ws.Range("C77").Value = 768
$ws.Range("D77").Value = 3448
$ws.Range("E77").Value = 9600
$ws.Range("G77").Value = 10
$ws.Range("H77").Value = 81
$ws.Range("B131").Value = 1710
$ws.Range("C131").Value = 21
$ws.Range("D131").Value = 889
$ws.Range("E131").Value = 816
$ws.Range("H131").Value = 5
$ws.Range("B132").Value = 1694
$ws.Range("C132").Value = 4
$ws.Range("D132").Value = 918
$ws.Range("E132").Value = 742
$ws.Range("H132").Value = 34
$ws.Range("B148").Value = 915
$ws.Range("C148").Value = 26
$ws.Range("D148").Value = 805
$ws.Range("E148").Value = 35
$ws.Range("H148").Value = 75
$ws.Range("D149").Value = 803
$ws.Range("E149").Value = 34
$ws.Range("H149").Value = 52
$ws.Range("B184").Value = 114
$ws.Range("C184").Value = 2
$ws.Range("E184").Value = 10
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 22
$ws.Range("E202").Value = 2
$ws.Range("B203").Value = 24
$ws.Range("D203").Value = 24
$ws.Range("E203").Value = 0

